# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values on Sheet1 for rows 2-39 to the newly regenerated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 7
    3  = 11
    4  = 5
    5  = 8
    6  = 7
    7  = 7
    8  = 6
    9  = 6
    10 = 7
    11 = 6
    12 = 1
    13 = 5
    14 = 12
    15 = 6
    16 = 8
    17 = 5
    18 = 10
    19 = 9
    20 = 10
    21 = 7
    22 = 6
    23 = 8
    24 = 3
    25 = 3
    26 = 5
    27 = 6
    28 = 6
    29 = 3
    30 = 5
    31 = 5
    32 = 6
    33 = 8
    34 = 13
    35 = 5
    36 = 6
    37 = 3
    38 = 0
    39 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
